# "Generate Report for Handoff"
# - Row for the 6100965e-3277-4e74-8ceb-b89abe4613f0 file is removed from every sheet
#   (it has been superseded / no longer part of the report).
# - The row for 18633056-ff34-44e2-8461-f8cb6b01ebaa is updated: status flips from
#   "Handed back: in sync with en-US" to "Ready for handoff" and the handoff
#   timestamps are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-03-19 16:48:24"

# drop all hyperlinks on the sheet, remove the obsolete row, then restore the
# hyperlink that is still needed
$ovLinkAddress = "https://github.com/OpenLocalizationTest/oltest/blob/321a715925ff888fed7ae077e120abccf0885bbf/e2e/18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$ov.Hyperlinks.Delete()
$ov.Range("A3:D3").EntireRow.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), $ovLinkAddress, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-19 16:48:14"

$zhLinkMd = "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/2b909e6746b9f4665f616d9dee9a80ec32acd39c/e2e/18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$zhLinkXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b9274bb20bef04eb0a0885ecd5ae2510223ad243/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf"
$zhLinkMdA = "https://github.com/OpenLocalizationTest/oltest/blob/321a715925ff888fed7ae077e120abccf0885bbf/e2e/18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$zhLinkXlfG = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/532b74cd7f400edd8f4f765e7060380762d88ab4/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf"

$zh.Hyperlinks.Delete()
$zh.Range("A3:L3").EntireRow.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhLinkMdA, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.md")
$zh.Hyperlinks.Add($zh.Range("D2"), $zhLinkXlf, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), $zhLinkMd, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.md")
$zh.Hyperlinks.Add($zh.Range("G2"), $zhLinkXlfG, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Ready for handoff"
# E2 uses the same shared timestamp text as Overview!D2

$deLinkMd = "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/b4aeaf9f3ee21d2171d9d36a8143a2ae168c271d/e2e/18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$deLinkXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e9c13ed9a961b8c5076d63cf9435422303951bd2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf"
$deLinkMdA = "https://github.com/OpenLocalizationTest/oltest/blob/321a715925ff888fed7ae077e120abccf0885bbf/e2e/18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$deLinkXlfG = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f772a49fcc6f11ce28d31189722e226908d1cbcc/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf"

$de.Hyperlinks.Delete()
$de.Range("A3:L3").EntireRow.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deLinkMdA, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.md")
$de.Hyperlinks.Add($de.Range("D2"), $deLinkXlf, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), $deLinkMd, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.md")
$de.Hyperlinks.Add($de.Range("G2"), $deLinkXlfG, [Type]::Missing, [Type]::Missing, "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf")

Write-Host "Handoff report regenerated"
